$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "编号"
$ws.Range("B1").Value = "姓名"
$ws.Range("C1").Value = "年龄"
$ws.Range("D1").Value = "注册日期"

# Data row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "asdf"
$ws.Range("C2").Value = 25
$ws.Range("D2").Value = 43101
$ws.Range("D2").NumberFormat = "mm-dd-yy"

# Data row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "cipchk"
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 43101

# Reuse D2's date style for D3 instead of minting a new cellXf
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("D4").Select()
